$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that LOOKS numeric (all digits) as genuine TEXT
# (shared string), matching how the source data was originally stored
# (RUT / DV / SECCION columns are text, not numbers). We route the value
# through a TEXT() formula and then Copy / Paste-Values it back onto
# itself; that yields a plain text cell with no residual formula and no
# extra number-format style, exactly like the rest of the sheet.
function Set-TextDigits($range, $digits) {
    $range.Formula = "=TEXT(" + $digits + ",""0"")"
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
}

# --- Update existing row 2 (uppercase the previously lowercase test data) ---
Set-TextDigits $ws.Range("B2") "20357099"
Set-TextDigits $ws.Range("C2") "6"
$ws.Range("D2").Value = "MATIAS IGNACIO"
$ws.Range("E2").Value = "CEBALLOS VASQUEZ"
# F2 (SECCION = 1) and G2 (ASIGNATURA) keep their existing values.

# --- New row 3 ---
$ws.Range("A3").Value = "31-12-2024"
Set-TextDigits $ws.Range("B3") "21075353"
Set-TextDigits $ws.Range("C3") "2"
$ws.Range("D3").Value = "FLAVIO ALEXANDER"
$ws.Range("E3").Value = "JARA LABRIN"
Set-TextDigits $ws.Range("F3") "1"
$ws.Range("G3").Value = "FÍSICA MECANICA / 3"

# --- New row 4 ---
$ws.Range("A4").Value = "31-12-2024"
Set-TextDigits $ws.Range("B4") "21075353"
Set-TextDigits $ws.Range("C4") "2"
$ws.Range("D4").Value = "ANA DORA"
$ws.Range("E4").Value = "LABRIN ESPINOZA"
Set-TextDigits $ws.Range("F4") "1"
$ws.Range("G4").Value = "FÍSICA MECANICA / 3"
